# edit.ps1 - apply the "added makefile in doc" change to UserManual.docx
#
# Summary of the edit (per the OOXML diff):
#   1. The list item that used to read "UML" + ".png" now reads "Makefile".
#   2. A brand-new list item "UML" + ".png" (two runs) is inserted right
#      after it, before the "UserManual.docx" list item.
#   3. The "Output will be similar to:" sentence is split into three runs,
#      with a grammar-check bracket (<w:proofErr w:type="gramStart"/> /
#      <w:proofErr w:type="gramEnd"/>) around "similar to".
#   4. Four words that look like spelling errors to Word's proofer
#      ("Userid", "mxotcyqjh", "fygurcvbi", "yfhpvjuhm") get wrapped with
#      <w:proofErr w:type="spellStart"/> / <w:proofErr w:type="spellEnd"/>.
#
# Strategy: whole-paragraph replacement via Range.InsertXML. Replacing a
# paragraph's full Range (the text plus its trailing paragraph mark) with a
# <w:p> fragment that carries an explicit <w:pPr> reliably substitutes that
# paragraph in place (preserving neighbours). New paragraphs are first
# created as plain-text placeholders via Range.InsertBefore (which
# correctly inherits the surrounding list formatting) and are then given
# their final run layout the same way.
#
# NOTE: this interpreter's argument binding is fragile for user-defined
# functions -- a parenthesised/concatenated expression passed straight as
# an argument silently turns into an empty value. Every value handed to a
# function below is therefore pre-computed into its own variable first.

$d = $word.ActiveDocument

function Set-RangeXml {
    param($Rng, [string]$InnerXml)
    $pkg = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $InnerXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $rv = $Rng.InsertXML($pkg)
}

function Find-ParagraphByText {
    param([string]$Like)
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text -like $Like) {
            return $p
        }
    }
    return $null
}

# ---------------------------------------------------------------------
# 1 & 2: "UML.png" list item -> "Makefile", plus a fresh "UML.png" item
#         right before "UserManual.docx".
# ---------------------------------------------------------------------

$listPPr = '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/></w:pPr>'

# Insert a placeholder paragraph before "UserManual.docx"; InsertBefore on
# the paragraph Range correctly inherits the ListParagraph/numPr formatting.
$userManualPara = Find-ParagraphByText "*UserManual.docx*"
$userManualPara.Range.InsertBefore("UMLPLACEHOLDER`r")

# Turn the placeholder into the real two-run "UML" / ".png" paragraph.
$placeholderPara = Find-ParagraphByText "UMLPLACEHOLDER`r"
$placeholderRng = $placeholderPara.Range
$umlPngInner = $listPPr + "<w:r><w:t>UML</w:t></w:r><w:r><w:t>.png</w:t></w:r></w:p>"
$umlPngInner = "<w:p>" + $umlPngInner
Set-RangeXml $placeholderRng $umlPngInner

# Turn the original "UML" + ".png" paragraph into "Makefile".
$umlPara = Find-ParagraphByText "UML.png`r"
$umlRng = $umlPara.Range
$makefileInner = $listPPr + "<w:r><w:t>Makefile</w:t></w:r></w:p>"
$makefileInner = "<w:p>" + $makefileInner
Set-RangeXml $umlRng $makefileInner

# ---------------------------------------------------------------------
# 3: "Output will be similar to:" gains a grammar-check proofErr bracket.
# ---------------------------------------------------------------------

$outputPara = Find-ParagraphByText "*Output will be similar to:*"
$outputRng = $outputPara.Range
$outputInner = '<w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:ind w:left="360"/></w:pPr>' +
    '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Output:</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> All output goes to the console. </w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">Output will be </w:t></w:r>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:t>similar to</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '<w:r><w:t>:</w:t></w:r>' +
    '</w:p>'
Set-RangeXml $outputRng $outputInner

# ---------------------------------------------------------------------
# 4: spellStart/spellEnd proofErr brackets around four flagged words.
# ---------------------------------------------------------------------

$useridPara = Find-ParagraphByText "*Userid*"
$useridRng = $useridPara.Range
$useridInner = '<w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:ind w:left="720"/></w:pPr>' +
    '<w:proofErr w:type="spellStart"/><w:r><w:t>Userid</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:tab/></w:r><w:r><w:tab/></w:r>' +
    '<w:r><w:t>Password</w:t></w:r>' +
    '<w:r><w:tab/></w:r>' +
    '<w:r><w:t>Resul</w:t></w:r><w:r><w:t>t</w:t></w:r>' +
    '</w:p>'
Set-RangeXml $useridRng $useridInner

$smithPara = Find-ParagraphByText "*mxotcyqjh*"
$smithRng = $smithPara.Range
$smithInner = '<w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:ind w:left="720"/></w:pPr>' +
    '<w:r><w:t>SMITH</w:t></w:r>' +
    '<w:r><w:tab/></w:r><w:r><w:tab/></w:r>' +
    '<w:proofErr w:type="spellStart"/><w:r><w:t>mxotcyqjh</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:tab/></w:r>' +
    '<w:r><w:t>match</w:t></w:r>' +
    '</w:p>'
Set-RangeXml $smithRng $smithInner

$johnsonPara = Find-ParagraphByText "*fygurcvbi*"
$johnsonRng = $johnsonPara.Range
$johnsonInner = '<w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:ind w:left="720"/></w:pPr>' +
    '<w:r><w:t>JOHNSON</w:t></w:r>' +
    '<w:r><w:tab/></w:r>' +
    '<w:proofErr w:type="spellStart"/><w:r><w:t>fygurcvbi</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:tab/></w:r>' +
    '<w:r><w:t>match</w:t></w:r>' +
    '</w:p>'
Set-RangeXml $johnsonRng $johnsonInner

$williamsPara = Find-ParagraphByText "*yfhpvjuhm*"
$williamsRng = $williamsPara.Range
$williamsInner = '<w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:ind w:left="720"/></w:pPr>' +
    '<w:r><w:t>WILLIAMS</w:t></w:r>' +
    '<w:r><w:tab/></w:r>' +
    '<w:proofErr w:type="spellStart"/><w:r><w:t>yfhpvjuhm</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:tab/></w:r>' +
    '<w:r><w:t>match</w:t></w:r>' +
    '</w:p>'
Set-RangeXml $williamsRng $williamsInner

Write-Output "edit complete"
